# chore: update Sheets via scheduled runner
# Refreshes cached market-board price/profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) for specific Leve rows across the per-job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2240.0833
$ws.Range("I19").Value = 1541.6666
$ws.Range("J19").Value = 2472.889
$ws.Range("K19").Value = 1541.6666
$ws.Range("L19").Value = 2472.889
$ws.Range("M19").Value = -1366.6666
$ws.Range("N19").Value = -2822.889
$ws.Range("H21").Value = 9994
$ws.Range("I21").Value = 9994
$ws.Range("K21").Value = 9994
$ws.Range("M21").Value = -9526
$ws.Range("H23").Value = 9994
$ws.Range("I23").Value = 9994
$ws.Range("K23").Value = 9994
$ws.Range("M23").Value = -9760
$ws.Range("H32").Value = 8033.0415
$ws.Range("J32").Value = 7805.263
$ws.Range("L32").Value = 7805.263
$ws.Range("N32").Value = -8457.262999999999
$ws.Range("H80").Value = 669.0476
$ws.Range("I80").Value = 566.0714
$ws.Range("J80").Value = 875
$ws.Range("K80").Value = 1698.2142
$ws.Range("L80").Value = 2625
$ws.Range("M80").Value = -700.2142000000001
$ws.Range("N80").Value = -4621
$ws.Range("H83").Value = 669.0476
$ws.Range("I83").Value = 566.0714
$ws.Range("J83").Value = 875
$ws.Range("K83").Value = 5094.6426
$ws.Range("L83").Value = 7875
$ws.Range("M83").Value = -102.6426000000001
$ws.Range("N83").Value = -17859
$ws.Range("H106").Value = 4352
$ws.Range("I106").Value = 3802.6667
$ws.Range("K106").Value = 3802.6667
$ws.Range("M106").Value = -3171.6667
$ws.Range("H116").Value = 4136.25
$ws.Range("I116").Value = 3560
$ws.Range("K116").Value = 3560
$ws.Range("M116").Value = -118
$ws.Range("H132").Value = 5274.28
$ws.Range("I132").Value = 5667.636
$ws.Range("J132").Value = 2389.6667
$ws.Range("K132").Value = 17002.908
$ws.Range("L132").Value = 7169.000100000001
$ws.Range("M132").Value = -14472.908
$ws.Range("N132").Value = -12229.0001
$ws.Range("H137").Value = 1379.3334
$ws.Range("I137").Value = 1211.7
$ws.Range("J137").Value = 1714.6
$ws.Range("K137").Value = 3635.1
$ws.Range("L137").Value = 5143.799999999999
$ws.Range("M137").Value = -1085.1
$ws.Range("N137").Value = -10243.8
$ws.Range("H138").Value = 2813.843
$ws.Range("I138").Value = 1826.64
$ws.Range("J138").Value = 3763.077
$ws.Range("K138").Value = 5479.92
$ws.Range("L138").Value = 11289.231
$ws.Range("M138").Value = -339.9200000000001
$ws.Range("N138").Value = -21569.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4488.5273
$ws.Range("I32").Value = 4386.463
$ws.Range("K32").Value = 4386.463
$ws.Range("M32").Value = -4099.463
$ws.Range("H45").Value = 107202.7
$ws.Range("I45").Value = 255555.88
$ws.Range("K45").Value = 255555.88
$ws.Range("M45").Value = -255178.88
$ws.Range("H61").Value = 13520662
$ws.Range("I61").Value = 16134911
$ws.Range("K61").Value = 16134911
$ws.Range("M61").Value = -16134699
$ws.Range("H110").Value = 5906.533
$ws.Range("I110").Value = 4327.1816
$ws.Range("K110").Value = 4327.1816
$ws.Range("M110").Value = -2282.1816
$ws.Range("H122").Value = 2668.375
$ws.Range("I122").Value = 2299
$ws.Range("K122").Value = 6897
$ws.Range("M122").Value = -4447
$ws.Range("H132").Value = 3760.2896
$ws.Range("I132").Value = 3108.818
$ws.Range("K132").Value = 9326.454000000002
$ws.Range("M132").Value = -6796.454000000002
$ws.Range("H136").Value = 13520662
$ws.Range("I136").Value = 16134911
$ws.Range("K136").Value = 48404733
$ws.Range("M136").Value = -48402183

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1390
$ws.Range("I105").Value = 1390
$ws.Range("K105").Value = 1390
$ws.Range("M105").Value = 357
$ws.Range("H119").Value = 32000
$ws.Range("J119").Value = 32000
$ws.Range("L119").Value = 32000
$ws.Range("N119").Value = -41676
$ws.Range("H135").Value = 95000
$ws.Range("J135").Value = 95000
$ws.Range("L135").Value = 95000
$ws.Range("N135").Value = -105140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2999.6667
$ws.Range("I4").Value = 2333
$ws.Range("J4").Value = 3666.3333
$ws.Range("K4").Value = 2333
$ws.Range("L4").Value = 3666.3333
$ws.Range("M4").Value = -2221
$ws.Range("N4").Value = -3890.3333
$ws.Range("H15").Value = 6148.5
$ws.Range("I15").Value = 2497.5
$ws.Range("K15").Value = 2497.5
$ws.Range("M15").Value = -2327.5
$ws.Range("H107").Value = 466.55554
$ws.Range("I107").Value = 385.2857
$ws.Range("K107").Value = 385.2857
$ws.Range("M107").Value = 1534.7143
$ws.Range("H132").Value = 3990.0557
$ws.Range("I132").Value = 3430.1428
$ws.Range("K132").Value = 10290.4284
$ws.Range("M132").Value = -7760.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 499.25
$ws.Range("J68").Value = 499.5
$ws.Range("L68").Value = 1498.5
$ws.Range("N68").Value = -3120.5
$ws.Range("H71").Value = 499.25
$ws.Range("J71").Value = 499.5
$ws.Range("L71").Value = 4495.5
$ws.Range("N71").Value = -12607.5
$ws.Range("H111").Value = 7999.5
$ws.Range("J111").Value = 7999.5
$ws.Range("L111").Value = 23998.5
$ws.Range("N111").Value = -30132.5
$ws.Range("H132").Value = 31251186
$ws.Range("J132").Value = 1480
$ws.Range("L132").Value = 13320
$ws.Range("N132").Value = -18380

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 405000600
$ws.Range("I11").Value = 506250000
$ws.Range("J11").Value = 3000
$ws.Range("K11").Value = 506250000
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = -506249861
$ws.Range("N11").Value = -3278
$ws.Range("H97").Value = 2874.375
$ws.Range("J97").Value = 7098.875
$ws.Range("L97").Value = 7098.875
$ws.Range("N97").Value = -8090.875
$ws.Range("H132").Value = 2756.3794
$ws.Range("I132").Value = 2355.5264
$ws.Range("J132").Value = 3518
$ws.Range("K132").Value = 7066.5792
$ws.Range("L132").Value = 10554
$ws.Range("M132").Value = -4536.5792
$ws.Range("N132").Value = -15614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3906.8147
$ws.Range("I16").Value = 3393.0667
$ws.Range("J16").Value = 4549
$ws.Range("K16").Value = 3393.0667
$ws.Range("L16").Value = 4549
$ws.Range("M16").Value = -3223.0667
$ws.Range("N16").Value = -4889
$ws.Range("H22").Value = 1845.1818
$ws.Range("J22").Value = 2124.625
$ws.Range("L22").Value = 2124.625
$ws.Range("N22").Value = -2714.625
$ws.Range("H27").Value = 1845.1818
$ws.Range("J27").Value = 2124.625
$ws.Range("L27").Value = 2124.625
$ws.Range("N27").Value = -2338.625
$ws.Range("H46").Value = 15209.125
$ws.Range("I46").Value = 3096.2856
$ws.Range("K46").Value = 3096.2856
$ws.Range("M46").Value = -2908.2856
$ws.Range("H123").Value = 68000
$ws.Range("J123").Value = 68000
$ws.Range("L123").Value = 68000
$ws.Range("N123").Value = -77800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2939.1428
$ws.Range("J2").Value = 2929.1667
$ws.Range("L2").Value = 2929.1667
$ws.Range("N2").Value = -3153.1667
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H107").Value = 1626.4615
$ws.Range("I107").Value = 974.2857
$ws.Range("K107").Value = 2922.8571
$ws.Range("M107").Value = -1002.8571
$ws.Range("H132").Value = 4347.14
$ws.Range("I132").Value = 3741.3948
$ws.Range("K132").Value = 11224.1844
$ws.Range("M132").Value = -8694.1844
